$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.127.59"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.048.44"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'248.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").Value = "'0.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("D7").Value = "'59.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("E10").Value = "  -2.42%  "

$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").Value = "'15.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.45%  "

$ws.Range("D13").Value = "2.346.19"
$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("D14").Value = "'0.840"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "

$ws.Range("E15").Value = "  +8.06%  "

$ws.Range("D16").Value = "2.057.54"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "'18.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +23.94%  "

$ws.Range("D18").Value = "37.122.62"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "'74.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").Value = "'5.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "'236.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'2.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").Value = "'169.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("E26").Value = "  +7.76%  "

$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  +5.90%  "

$ws.Range("E31").Value = "  +3.05%  "

$ws.Range("D32").Value = "'0.0624"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("D33").Value = "'4.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "

$ws.Range("D34").Value = "'0.0897"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("E38").Value = "  -2.71%  "

$ws.Range("D39").Value = "'1.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "

$ws.Range("D40").Value = "'3.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.56%  "

$ws.Range("D41").Value = "'5.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.47%  "

$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").Value = "'17.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.18%  "

$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'95.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("D48").Value = "1.281.56"
$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("D49").Value = "'6.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "

$ws.Range("D50").Value = "2.235.72"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -19.09%  "
